$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "29.003.72"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.65%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.830.94"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.42%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.9995"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "241.48"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.30%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.6261"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -5.54%  "
$ws.Range("E7").Value = "  +0.07%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "45.25"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +2.03%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07594"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +2.14%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.2915"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.09%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "22.71"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.18%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.07746"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.15%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.831.40"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.51%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "4.954"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.39%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.6626"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.46%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "82.32"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.18%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.000009468"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +9.44%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "5.980"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -2.84%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "28.989.74"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.73%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "223.98"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("E21").Value = "  -1.62%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "1.0000"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "7.201"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("E24").Value = "  +0.07%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "158.92"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.19%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "8.408"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.41%  "
$ws.Range("E27").Value = "  -3.35%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "17.82"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.27%  "
$ws.Range("E29").Value = "  -1.28%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.052"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.91%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.022"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.85%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.191"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("E33").Value = "  -2.63%  "
$ws.Range("E34").Value = "  -1.93%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.7363"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("E36").Value = "  -0.86%  "
$ws.Range("E37").Value = "  +1.79%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.253.02"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -4.75%  "
$ws.Range("E39").Value = "  +0.77%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.01787"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.61%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "6.206"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -3.18%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.8872"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.40%  "
$ws.Range("E43").Value = "  +0.20%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "101.43"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.48%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.974.83"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.98%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "64.52"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.90%  "
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("E48").Value = "  +0.28%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.3974"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.03%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "8.836"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("E51").Value = "  -1.67%  "
